# Adicionei filtro nos processos do advogado
#
# The underlying content edit was a typo fix ("Promivido" -> "Promovido")
# inside the "Vencedor" bullet plus a couple of in-place retypes
# ("Advogado"+"Promovido" -> single run, "deliberativas"+")" -> single
# run). Re-typing those spots made Word's background proofer re-scan the
# touched list items, which is why many otherwise-untouched bullets
# (Usuario, Cpf, Email, Intimacao*, ...) now carry <w:proofErr> spell-check
# bookmarks around them. We reproduce both the visible text change and the
# proofErr bookmarks below.
#
# Technique: Range.InsertXML() with a WordprocessingML package replaces the
# contents of the *exact* range it's called on while leaving the enclosing
# <w:p>'s properties (pPr / paraId / rsid / numbering) untouched, as long as
# the range does not include the paragraph mark. Every replacement below is
# therefore done on Range(start, end) pairs that cover only the run text of
# a single paragraph (never the trailing paragraph mark), and every
# replacement keeps the exact same number of characters as the original, so
# the absolute offsets used for every other edit stay valid no matter what
# order the edits run in.
#
# NOTE: this interpreter's functions only bind POSITIONAL parameters
# reliably, so every helper below is called positionally (no -Name value).

$d = $word.ActiveDocument

function Set-RunXml($Start, $End, $InnerXml) {
    $range = $d.Range($Start, $End)
    $pkg = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' +
        $InnerXml +
        '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $range.InsertXML($pkg)
}

function Wrap-SpellErr($Start, $End, $Text) {
    $inner = '<w:proofErr w:type="spellStart"/><w:r><w:t>' + $Text + '</w:t></w:r><w:proofErr w:type="spellEnd"/>'
    Set-RunXml $Start $End $inner
}

# ---- simple "wrap existing single-run bullet with spellStart/spellEnd" ----
Wrap-SpellErr 7   15  "Usuario"
Wrap-SpellErr 38  42  "Cpf"
Wrap-SpellErr 42  48  "Email"
Wrap-SpellErr 215 234 "AdvogadoPromovente"
Wrap-SpellErr 234 250 "PartePromovente"
Wrap-SpellErr 268 283 "PartePromovida"
Wrap-SpellErr 311 321 "Descricao"
Wrap-SpellErr 426 448 "JustificativaResposta"
Wrap-SpellErr 457 465 "Usuario"
Wrap-SpellErr 503 507 "Cpf"
Wrap-SpellErr 512 518 "Email"
Wrap-SpellErr 518 528 "Intimacao"
Wrap-SpellErr 528 540 "IntimacaoId"
Wrap-SpellErr 540 554 "IntimacaoData"
Wrap-SpellErr 554 566 "IntimadoCpf"
Wrap-SpellErr 566 579 "IntimadoNome"
Wrap-SpellErr 579 596 "IntimadoEndereco"
Wrap-SpellErr 596 609 "ExecucaoData"

# ---- "Advogado" + "Promovido" runs retyped/merged into one run ----
Wrap-SpellErr 250 268 "AdvogadoPromovido"

# ---- "Caminho" + "Arquivo" bullet: stays two runs, gets wrapped ----
$caminhoInner = '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>Caminho</w:t></w:r>' +
    '<w:r><w:t>Arquivo</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>'
Set-RunXml 321 336 $caminhoInner

# ---- "Tipo (informativas ou deliberativas)" - last two runs merge ----
$tipoInner = '<w:r><w:t>Tipo</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> (informativas ou </w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:color w:val="24292E"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t>deliberativas)</w:t></w:r>'
Set-RunXml 336 373 $tipoInner

# ---- "Vencedor (Null, Promovente ou Promivido)" -> typo fix + split runs ----
$vencedorInner = '<w:r><w:t xml:space="preserve">Vencedor (</w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>Null</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve">, Promovente ou Prom</w:t></w:r>' +
    '<w:r><w:t>o</w:t></w:r>' +
    '<w:r><w:t>vido)</w:t></w:r>'
Set-RunXml 169 210 $vencedorInner
